$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - first worksheet
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 416
$ws1.Range("F3").Value = 2525
$ws1.Range("F4").Value = 117

# Sheet "全部类型" (all types) - fourth worksheet
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 416
$ws4.Range("F7").Value = 2525
$ws4.Range("F8").Value = 117
